# Applies the "Ver3" corrections to the furniture sale test-specification sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update test-case texts / expected-result wording ---

# D23: "result" -> fuller description of the expected result column
$ws.Range("D23").Value = "result (Decimal - precision 2 - Round_Half_Up)"

# J24 / J29: the generic "Einfacher String" note is split into a
# total-specific and a day-specific variant
$ws.Range("J24").Value = "Einfacher String für total"
$ws.Range("J29").Value = "Einfacher String für day"

# D30,31,38,39,46,47,54,55: append the computed result "=0.15" to the
# Decimal("0.05") expression test cases
$dVal015 = 'Decimal("0.05")+Decimal("0.1")-(Decimal("0.05")*Decimal("0.1"))=0.15'
foreach ($r in 30,31,38,39,46,47,54,55) {
    $ws.Range("D$r").Value = $dVal015
}

# D62,63,70,71,78,79,86,87: append the computed result "=0.19" to the
# Decimal("0.1") expression test cases
$dVal019 = 'Decimal("0.1")+Decimal("0.1")-(Decimal("0.1")*Decimal("0.1"))=0.19'
foreach ($r in 62,63,70,71,78,79,86,87) {
    $ws.Range("D$r").Value = $dVal019
}

# D94,95,102,103: append the computed result "=0.28" to the
# Decimal("0.2") expression test cases
$dVal028 = 'Decimal("0.2")+Decimal("0.1")-(Decimal("0.2")*Decimal("0.1"))=0.28'
foreach ($r in 94,95,102,103) {
    $ws.Range("D$r").Value = $dVal028
}

# Remove stray leftover side-note cells in column I (rows 45-47)
$ws.Range("I45").ClearContents()
$ws.Range("I46").ClearContents()
$ws.Range("I47").ClearContents()

# --- Update the sheet view (scrolled/zoomed/selected differently) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
$win.Zoom = 115
$ws.Range("F86").Select()
